# Delete columns L:O (shifts dimension to A1:K5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1:O5").Delete() | Out-Null

# Update header labels for F1:K1
$ws.Range("F1").Value = "Fit time STDT"
$ws.Range("G1").Value = "Prediction time STDT"
$ws.Range("H1").Value = "Score STDT"
$ws.Range("I1").Value = "Fit time CondensedDT"
$ws.Range("J1").Value = "Prediction time CondensedDT"
$ws.Range("K1").Value = "Score CondensedDT"

# Row 2 (Iris)
$ws.Range("F2").Value = 0.01599884033203125
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0.02599978446960449
$ws.Range("J2").Value = 0.003999948501586914
$ws.Range("K2").Value = 0.3

# Row 3 (Wine)
$ws.Range("F3").Value = 0.1229982376098633
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.9166666666666666
$ws.Range("I3").Value = 0.04700374603271484
$ws.Range("J3").Value = 0.01999711990356445
$ws.Range("K3").Value = 0.3055555555555556

# Row 4 (Breast Cancer)
$ws.Range("F4").Value = 2.365621328353882
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.9298245614035088
$ws.Range("I4").Value = 0.1370000839233398
$ws.Range("J4").Value = 0.1309998035430908
$ws.Range("K4").Value = 0.5526315789473685

# Row 5 (Digits)
$ws.Range("F5").Value = 1.332000017166138
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.8694444444444445
$ws.Range("I5").Value = 1.445001840591431
$ws.Range("J5").Value = 1.197999477386475
$ws.Range("K5").Value = 0.2083333333333333
